$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.036339451647734222,
    -0.0099999991695227664,
    -0.0089999991793181522,
    0.28398587613580872,
    -0.0059999992042714112,
    -0.0059999991779662309,
    -0.019999999038734728,
    -0.019999999034784111,
    -0.0059999991697798905,
    -0.0059999991682460063,
    -0.0044999991829541841,
    -0.0059999991683001852,
    -0.0059999991702284206,
    0.046712874415856476,
    -0.0059999991684254184,
    -0.0059999991653496565,
    -0.0059999991608945535,
    -0.014090005254781879,
    -0.0089999991790987721,
    -0.008999999172178974,
    -0.0089999991711744443,
    -0.0089999991704692306,
    -0.0089999991728095807,
    -0.041999998835126284,
    -0.041999998829107099,
    -0.0059999991747794468,
    -0.0059999991704700051,
    -0.0059999991517036833,
    -0.011999999079678858,
    -0.019999998994850277,
    -0.014999999034323253,
    -0.020999998973616485,
    -0.014197333089365216
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
